# "@dev finish 5 interfaces"
# Update the worklist: rename the "get user info" task to the asyncstorage
# version, and mark 5 tasks as completed ("完成").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 task description changed from the axios-based task to the
# asyncstorage-based one. Set this first so the new shared string for
# "asyncstorage 获取用户信息" is registered before the "完成" string.
$ws.Range("B3").Value = "asyncstorage 获取用户信息"

# Mark 5 tasks as finished (status column D), in row order.
$ws.Range("D2").Value = "完成"
$ws.Range("D3").Value = "完成"
$ws.Range("D5").Value = "完成"
$ws.Range("D6").Value = "完成"
$ws.Range("D9").Value = "完成"

# Update the active cell selection on the sheet to reflect where the
# author left off editing.
$ws.Range("C21").Select() | Out-Null
